$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.917.29"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "3.331.54"
$ws.Range("E3").Value = "  +1.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.37"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.28"
$ws.Range("E6").Value = "  +2.03%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("E8").Value = "  +1.38%  "

$ws.Range("D9").Value = "3.328.02"
$ws.Range("E9").Value = "  +1.49%  "

$ws.Range("E10").Value = "  +5.54%  "

$ws.Range("E11").Value = "  +1.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.47"
$ws.Range("E12").Value = "  +4.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000274"
$ws.Range("E13").Value = "  +2.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "705.74"
$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").Value = "3.866.44"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("D17").Value = "67.946.97"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("D19").Value = "3.329.37"
$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("E20").Value = "  +1.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.08"
$ws.Range("E21").Value = "  +3.17%  "

$ws.Range("E22").Value = "  +0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.39"
$ws.Range("E23").Value = "  +4.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.07"
$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.93"
$ws.Range("E25").Value = "  +0.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.92"
$ws.Range("E26").Value = "  +1.01%  "

$ws.Range("E27").Value = "  +0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.61"
$ws.Range("E28").Value = "  +4.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.10"
$ws.Range("E29").Value = "  -1.80%  "

$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("E31").Value = "  +6.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "569.37"
$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("E34").Value = "  +3.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.33"
$ws.Range("E35").Value = "  +3.63%  "

$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.44"
$ws.Range("E36").Value = "  +4.67%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").Value = "3.684.29"
$ws.Range("E38").Value = "  -4.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.76"
$ws.Range("E39").Value = "  +9.87%  "

$ws.Range("E40").Value = "  +3.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  +3.20%  "

$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.34"
$ws.Range("E44").Value = "  -2.48%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.337"
$ws.Range("E45").Value = "  +3.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0407"
$ws.Range("E46").Value = "  +0.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  +5.71%  "

$ws.Range("E48").Value = "  +1.31%  "

$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("E50").Value = "  -2.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.29"
$ws.Range("E51").Value = "  +0.78%  "
